# Auto-generated edit script: update Halicarnassus_Profits price/profit data
# across sheets ALC, ARM, CRP, CUL, GSM, LTW, WVR (BSM untouched).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = ""
$ws.Range("H40").Value = 6748.593
$ws.Range("I40").Value = 4499.7
$ws.Range("J40").Value = 8071.4707
$ws.Range("K40").Value = 4499.7
$ws.Range("L40").Value = 8071.4707
$ws.Range("M40").Value = -4324.7
$ws.Range("N40").Value = -8421.4707
$ws.Range("H94").Value = 2552.9167
$ws.Range("I94").Value = 2552.9167
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2552.9167
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2101.9167
$ws.Range("N94").Value = ""
$ws.Range("H112").Value = 4642.857
$ws.Range("J112").Value = 5000
$ws.Range("L112").Value = 15000
$ws.Range("N112").Value = -17216
$ws.Range("H113").Value = 5741.5835
$ws.Range("J113").Value = 6189.4
$ws.Range("L113").Value = 6189.4
$ws.Range("N113").Value = -12697.4
$ws.Range("H137").Value = 3091.1
$ws.Range("I137").Value = 1682.5
$ws.Range("K137").Value = 5047.5
$ws.Range("M137").Value = -2497.5
$ws.Range("H138").Value = 4633.8945
$ws.Range("I138").Value = 4318.2
$ws.Range("J138").Value = 4746.643
$ws.Range("K138").Value = 12954.6
$ws.Range("L138").Value = 14239.929
$ws.Range("M138").Value = -7814.599999999999
$ws.Range("N138").Value = -24519.929
$ws.Range("H141").Value = 6950
$ws.Range("I141").Value = 6950
$ws.Range("K141").Value = 20850
$ws.Range("M141").Value = -15670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 18501500
$ws.Range("I6").Value = 17859286
$ws.Range("K6").Value = 17859286
$ws.Range("M6").Value = -17859113
$ws.Range("H32").Value = 2383610
$ws.Range("I32").Value = 2520.7715
$ws.Range("J32").Value = 14289056
$ws.Range("K32").Value = 2520.7715
$ws.Range("L32").Value = 14289056
$ws.Range("M32").Value = -2233.7715
$ws.Range("N32").Value = -14289630
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H122").Value = 2305.25
$ws.Range("I122").Value = 2999.5
$ws.Range("J122").Value = 1611
$ws.Range("K122").Value = 8998.5
$ws.Range("L122").Value = 4833
$ws.Range("M122").Value = -6548.5
$ws.Range("N122").Value = -9733

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2065.2307
$ws.Range("J68").Value = 2452.3333
$ws.Range("L68").Value = 7356.999899999999
$ws.Range("N68").Value = -8978.999899999999
$ws.Range("H71").Value = 2065.2307
$ws.Range("J71").Value = 2452.3333
$ws.Range("L71").Value = 22070.9997
$ws.Range("N71").Value = -30182.9997
$ws.Range("H86").Value = 430.375
$ws.Range("I86").Value = 369.4
$ws.Range("J86").Value = 532
$ws.Range("K86").Value = 1108.2
$ws.Range("L86").Value = 1596
$ws.Range("M86").Value = 77.80000000000018
$ws.Range("N86").Value = -3968
$ws.Range("H89").Value = 430.375
$ws.Range("I89").Value = 369.4
$ws.Range("J89").Value = 532
$ws.Range("K89").Value = 3324.6
$ws.Range("L89").Value = 4788
$ws.Range("M89").Value = 2603.4
$ws.Range("N89").Value = -16644
$ws.Range("H107").Value = 616.5
$ws.Range("I107").Value = 690.3333
$ws.Range("J107").Value = 395
$ws.Range("K107").Value = 2070.9999
$ws.Range("L107").Value = 1185
$ws.Range("M107").Value = -150.9998999999998
$ws.Range("N107").Value = -5025
$ws.Range("H141").Value = 1641.8
$ws.Range("I141").Value = 1052.25
$ws.Range("K141").Value = 3156.75
$ws.Range("M141").Value = 2023.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 34000
$ws.Range("J15").Value = 34000
$ws.Range("L15").Value = 34000
$ws.Range("N15").Value = -34576
$ws.Range("H80").Value = 5264.8
$ws.Range("I80").Value = 6929.6
$ws.Range("J80").Value = 3600
$ws.Range("K80").Value = 6929.6
$ws.Range("L80").Value = 3600
$ws.Range("M80").Value = -5931.6
$ws.Range("N80").Value = -5596
$ws.Range("H81").Value = 34000
$ws.Range("J81").Value = 34000
$ws.Range("L81").Value = 34000
$ws.Range("N81").Value = -35996
$ws.Range("H83").Value = 5264.8
$ws.Range("I83").Value = 6929.6
$ws.Range("J83").Value = 3600
$ws.Range("K83").Value = 34648
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = -29656
$ws.Range("N83").Value = -27984
$ws.Range("H84").Value = 34000
$ws.Range("J84").Value = 34000
$ws.Range("L84").Value = 102000
$ws.Range("N84").Value = -111984
$ws.Range("H113").Value = 5047
$ws.Range("I113").Value = 2578.3333
$ws.Range("K113").Value = 2578.3333
$ws.Range("M113").Value = -408.3332999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3066.6667
$ws.Range("I22").Value = 2100
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 2100
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -1805
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 3066.6667
$ws.Range("I27").Value = 2100
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 2100
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -1993
$ws.Range("N27").Value = -5214
$ws.Range("H40").Value = 4994.75
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5272
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11470
$ws.Range("J62").Value = 12112.5
$ws.Range("L62").Value = 12112.5
$ws.Range("N62").Value = -13360.5
$ws.Range("H65").Value = 11470
$ws.Range("J65").Value = 12112.5
$ws.Range("L65").Value = 60562.5
$ws.Range("N65").Value = -66802.5
$ws.Range("H81").Value = 4426.8184
$ws.Range("I81").Value = 3369.5
$ws.Range("J81").Value = 15000
$ws.Range("K81").Value = 6739
$ws.Range("L81").Value = 30000
$ws.Range("M81").Value = -5678
$ws.Range("N81").Value = -32122
$ws.Range("H84").Value = 4426.8184
$ws.Range("I84").Value = 3369.5
$ws.Range("J84").Value = 15000
$ws.Range("K84").Value = 33695
$ws.Range("L84").Value = 150000
$ws.Range("M84").Value = -28391
$ws.Range("N84").Value = -160608
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H122").Value = 2666.4827
$ws.Range("I122").Value = 2705.2307
$ws.Range("K122").Value = 8115.6921
$ws.Range("M122").Value = -5665.6921

